$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NIT values for A2:A11 (replacing the previous, differently-sorted list)
$values = @(800057394, 900849952, 830074930, 800238184, 830071086, 900051936, 860530921, 900745947, 900323880, 900643491)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Rows 12-21 are no longer used -> clear their values (formatting/style stays as-is)
$ws.Range("A12:A21").ClearContents()

# A2:A11 previously carried the highlighted (yellow-fill) style; drop the
# highlight so the whole A2:A21 range shares the same plain bordered style
# that A12:A21 already used. Copy/PasteSpecial(Formats) reuses the existing
# style slot instead of inventing a new one.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A2:A11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# View: drop the frozen/scrolled topLeftCell and shrink the selection to
# the now-populated range
$ws.Range("A1").Select() | Out-Null
$ws.Range("A2:A11").Select() | Out-Null
